$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.782.38'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '2.536.80'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'317.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.08%  '
$ws.Range('D6').Value = "'95.43"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.25%  '
$ws.Range('D7').Value = "'0.580"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'0.530"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.05%  '
$ws.Range('D10').Value = "'36.27"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').Value = "'0.0811"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('D12').Value = "'7.60"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = '2.925.08'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.548.82'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = "'15.49"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.21%  '
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').Value = '42.719.14'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').Value = "'13.02"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('D20').Value = "'6.59"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').Value = '0.0₃0964'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = "'70.37"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.73%  '
$ws.Range('D23').Value = "'251.80"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('D24').Value = "'2.97"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('D25').Value = "'2.02"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range('D26').Value = "'26.86"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('E28').Value = '  +4.62%  '
$ws.Range('D29').Value = "'39.25"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.36%  '
$ws.Range('D30').Value = "'10.18"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.81%  '
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').Value = "'154.67"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.53%  '
$ws.Range('E33').Value = '  +1.11%  '
$ws.Range('D34').Value = "'18.98"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.80%  '
$ws.Range('D35').Value = "'3.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('D37').Value = "'2.64"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').Value = "'0.112"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.99%  '
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').Value = "'23.59"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.91%  '
$ws.Range('E41').Value = '  +10.69%  '
$ws.Range('D42').Value = "'3.80"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.61%  '
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('D44').Value = "'0.0302"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').Value = "'3.31"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.20%  '
$ws.Range('D46').Value = '2.019.25'
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('D47').Value = "'85.99"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('D49').Value = '2.779.25'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').Value = "'74.24"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.66%  '
$ws.Range('D51').Value = "'102.87"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.31%  '
